$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.348055371403689
$ws.Range("C2").Value = 0.2628968873353301
$ws.Range("D2").Value = 0.6557819579275304
$ws.Range("E2").Value = 0.2675972005155671
$ws.Range("G2").Value = 0.002470267621967948
$ws.Range("J2").Value = 0.1393798719812906
$ws.Range("N2").Value = 1.221051011966832
$ws.Range("O2").Value = 4.437661690615897

$ws.Range("B3").Value = 1.239695700604386
$ws.Range("C3").Value = 0.2349589806488552
$ws.Range("D3").Value = 0.6455821278322844
$ws.Range("E3").Value = 0.2625688216117723
$ws.Range("G3").Value = 0.002474190087930312
$ws.Range("J3").Value = 0.1360303736405513
$ws.Range("N3").Value = 1.236415440659222
$ws.Range("O3").Value = 4.42089757341185

$ws.Range("B4").Value = 1.17356139245976
$ws.Range("C4").Value = 0.2178260734181379
$ws.Range("D4").Value = 0.6396699616170451
$ws.Range("E4").Value = 0.2596309926111005
$ws.Range("G4").Value = 0.002476726338975978
$ws.Range("J4").Value = 0.1340563102826238
$ws.Range("N4").Value = 1.246394274039631
$ws.Range("O4").Value = 4.413634057831871

$ws.Range("B5").Value = 1.146712050569192
$ws.Range("C5").Value = 0.2108497353663097
$ws.Range("D5").Value = 0.6373487835328717
$ws.Range("E5").Value = 0.2584713667083136
$ws.Range("G5").Value = 0.002477792132379401
$ws.Range("J5").Value = 0.1332725640220005
$ws.Range("N5").Value = 1.250597730104886
$ws.Range("O5").Value = 4.411433926174595

$ws.Range("B6").Value = 1.142259851547351
$ws.Range("C6").Value = 0.2096916543624729
$ws.Range("D6").Value = 0.6369686712002363
$ws.Range("E6").Value = 0.2582810787505991
$ws.Range("G6").Value = 0.00247797105755181
$ws.Range("J6").Value = 0.133143672218452
$ws.Range("N6").Value = 1.251303984002259
$ws.Range("O6").Value = 4.411114427340436

$ws.Range("B7").Value = 1.173198883399721
$ws.Range("C7").Value = 0.2177319657565988
$ws.Range("D7").Value = 0.6396383008723774
$ws.Range("E7").Value = 0.2596152014746451
$ws.Range("G7").Value = 0.002476740581902594
$ws.Range("J7").Value = 0.134045656662714
$ws.Range("N7").Value = 1.246450408674395
$ws.Range("O7").Value = 4.413601312365302

$ws.Range("B8").Value = 1.31061047343303
$ws.Range("C8").Value = 0.2532596082251359
$ws.Range("D8").Value = 0.6521922441298784
$ws.Range("E8").Value = 0.2658323246812699
$ws.Range("G8").Value = 0.002471593617456525
$ws.Range("J8").Value = 0.1382077984765004
$ws.Range("N8").Value = 1.226235521116131
$ws.Range("O8").Value = 4.431251355968669

$ws.Range("B9").Value = 1.583228403589487
$ws.Range("C9").Value = 0.3230931164574145
$ws.Range("D9").Value = 0.6795977730112952
$ws.Range("E9").Value = 0.2792148934170413
$ws.Range("G9").Value = 0.00246251004796067
$ws.Range("J9").Value = 0.1470277948504943
$ws.Range("N9").Value = 1.190920491936879
$ws.Range("O9").Value = 4.490003301918989

$ws.Range("B10").Value = 1.785450617833419
$ws.Range("C10").Value = 0.3745006453180508
$ws.Range("D10").Value = 0.7014428757124165
$ws.Range("E10").Value = 0.2897795712128826
$ws.Range("G10").Value = 0.002456445181701357
$ws.Range("J10").Value = 0.1539142153309996
$ws.Range("N10").Value = 1.167614206818115
$ws.Range("O10").Value = 4.548037263180504

$ws.Range("B11").Value = 1.877868854367932
$ws.Range("C11").Value = 0.3979097483191367
$ws.Range("D11").Value = 0.7117547719169295
$ws.Range("E11").Value = 0.2947462956441456
$ws.Range("G11").Value = 0.002453816900889723
$ws.Range("J11").Value = 0.1571364530004047
$ws.Range("N11").Value = 1.157585348399216
$ws.Range("O11").Value = 4.577700959428455

$ws.Range("B12").Value = 1.912926323236377
$ws.Range("C12").Value = 0.4067774866508671
$ws.Range("D12").Value = 0.7157136175037238
$ws.Range("E12").Value = 0.2966502835726672
$ws.Range("G12").Value = 0.002452840320236222
$ws.Range("J12").Value = 0.1583695898897304
$ws.Range("N12").Value = 1.153870206761567
$ws.Range("O12").Value = 4.589405700966267

$ws.Range("B13").Value = 1.905373380618869
$ws.Range("C13").Value = 0.4048675200183993
$ws.Range("D13").Value = 0.7148586078744756
$ws.Range("E13").Value = 0.2962391928049044
$ws.Range("G13").Value = 0.002453049814384177
$ws.Range("J13").Value = 0.1581034349969457
$ws.Range("N13").Value = 1.154666654522764
$ws.Range("O13").Value = 4.586863855777153

$ws.Range("B14").Value = 1.880751840637799
$ws.Range("C14").Value = 0.3986392389283537
$ws.Range("D14").Value = 0.7120793865666712
$ws.Range("E14").Value = 0.2949024725394978
$ws.Range("G14").Value = 0.002453736183062458
$ws.Range("J14").Value = 0.1572376440916941
$ws.Range("N14").Value = 1.157278045870363
$ws.Range("O14").Value = 4.578654447186068

$ws.Range("B15").Value = 1.86567832177451
$ws.Range("C15").Value = 0.3948246511900493
$ws.Range("D15").Value = 0.7103840629942795
$ws.Range("E15").Value = 0.2940867159703089
$ws.Range("G15").Value = 0.002454159035057489
$ws.Range("J15").Value = 0.1567090099920421
$ws.Range("N15").Value = 1.158888355460903
$ws.Range("O15").Value = 4.573687454469507

$ws.Range("B16").Value = 1.779419394551383
$ws.Range("C16").Value = 0.372971264958835
$ws.Range("D16").Value = 0.7007765131636745
$ws.Range("E16").Value = 0.2894582261155492
$ws.Range("G16").Value = 0.002456619566971519
$ws.Range("J16").Value = 0.1537054424690467
$ws.Range("N16").Value = 1.168281163380811
$ws.Range("O16").Value = 4.546164552335256

$ws.Range("B17").Value = 1.726611087497417
$ws.Range("C17").Value = 0.3595708334009942
$ws.Range("D17").Value = 0.6949785745345878
$ws.Range("E17").Value = 0.2866600346073582
$ws.Range("G17").Value = 0.002458162419221606
$ws.Range("J17").Value = 0.1518858383935964
$ws.Range("N17").Value = 1.174190270732403
$ws.Range("O17").Value = 4.530117784538447

$ws.Range("B18").Value = 1.696277310968981
$ws.Range("C18").Value = 0.3518654870693467
$ws.Range("D18").Value = 0.6916789918421671
$ws.Range("E18").Value = 0.2850657195393111
$ws.Range("G18").Value = 0.002459062130934236
$ws.Range("J18").Value = 0.1508476770939779
$ws.Range("N18").Value = 1.177642985256149
$ws.Range("O18").Value = 4.521195162496326

$ws.Range("B19").Value = 1.686013732876518
$ws.Range("C19").Value = 0.3492569769752549
$ws.Range("D19").Value = 0.6905678587717716
$ws.Range("E19").Value = 0.2845285083155886
$ws.Range("G19").Value = 0.002459368873948629
$ws.Range("J19").Value = 0.1504976188587221
$ws.Range("N19").Value = 1.178821277525699
$ws.Range("O19").Value = 4.518226779447218

$ws.Range("B20").Value = 1.732228469539507
$ws.Range("C20").Value = 0.3609971029122789
$ws.Range("D20").Value = 0.6955921271959369
$ws.Range("E20").Value = 0.2869563403678015
$ws.Range("G20").Value = 0.002457996907388016
$ws.Range("J20").Value = 0.1520786657368944
$ws.Range("N20").Value = 1.173555651402594
$ws.Range("O20").Value = 4.531794194025792

$ws.Range("B21").Value = 1.8879821415577
$ws.Range("C21").Value = 0.4004685500231631
$ws.Range("D21").Value = 0.7128942460468295
$ws.Range("E21").Value = 0.2952944695192627
$ws.Range("G21").Value = 0.002453534073356266
$ws.Range("J21").Value = 0.1574915960081853
$ws.Range("N21").Value = 1.156508774939482
$ws.Range("O21").Value = 4.581052927118492

$ws.Range("B22").Value = 1.990129631499485
$ws.Range("C22").Value = 0.4262841098836248
$ws.Range("D22").Value = 0.7245167205864789
$ws.Range("E22").Value = 0.3008791654036429
$ws.Range("G22").Value = 0.002450726259590978
$ws.Range("J22").Value = 0.1611047455948693
$ws.Range("N22").Value = 1.145848999192367
$ws.Range("O22").Value = 4.615997176734879

$ws.Range("B23").Value = 1.935579452685261
$ws.Range("C23").Value = 0.4125042046890712
$ws.Range("D23").Value = 0.7182847699868944
$ws.Range("E23").Value = 0.2978861080814283
$ws.Range("G23").Value = 0.002452214909597614
$ws.Range("J23").Value = 0.1591694108049069
$ws.Range("N23").Value = 1.151494233984074
$ws.Range("O23").Value = 4.597094268134981

$ws.Range("B24").Value = 1.729688768767289
$ws.Range("C24").Value = 0.3603522903795238
$ws.Range("D24").Value = 0.6953146349860617
$ws.Range("E24").Value = 0.2868223356787993
$ws.Range("G24").Value = 0.002458071695660588
$ws.Range("J24").Value = 0.1519914637279527
$ws.Range("N24").Value = 1.17384238996334
$ws.Range("O24").Value = 4.531035346039289

$ws.Range("B25").Value = 1.509139683936155
$ws.Range("C25").Value = 0.3041836717670208
$ws.Range("D25").Value = 0.6718843213220964
$ws.Range("E25").Value = 0.2754663971041396
$ws.Range("G25").Value = 0.002464859998370672
$ws.Range("J25").Value = 0.1445707685890554
$ws.Range("N25").Value = 1.200010778155196
$ws.Range("O25").Value = 4.471509287091209
